$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Header "  Nombre de campaña" -> "Nombre de campaña" (leading spaces removed)
$ws.Range("B1").Value = "Nombre de campaña"

# Header "Revenue" -> "Ingresos"
$ws.Range("F1").Value = "Ingresos"

# Campaign type "de la experiencia del usuario (CEIP)" -> "Experiencia del cliente"
$ws.Range("D5").Value = "Experiencia del cliente"
$ws.Range("D8").Value = "Experiencia del cliente"
$ws.Range("D13").Value = "Experiencia del cliente"
